$d = $word.ActiveDocument

# Locate the "_GoBack" bookmark that sits right after the paragraph whose
# text is exactly "12345" (immediately before the bookmark markers).
$bm = $d.Bookmarks.Item("_GoBack")
$splitPos = $bm.Start

# Insert the new text "5465654" immediately before the bookmark, using the
# bookmark's own Range so the bookmark correctly re-anchors after the
# inserted text instead of swallowing it.
$bmRange = $bm.Range
$bmRange.InsertBefore("5465654")

# Now split the paragraph at the boundary between the original "12345" run
# and the newly inserted "5465654" text, turning them into two separate
# paragraphs (the bookmark travels with the second one).
$rSplit = $d.Range($splitPos, $splitPos)
$rSplit.InsertParagraphAfter()
